# v1: undo/redo apply color
#
# Task #33 "undo, redo coloring a section on the image" is finished -
# move it from the Active sheet to the Inactive sheet (mark Done, stamp
# today's completion date), and log two new follow-on Active tasks that
# came out of doing that work (#52, #53). Bump the tracked Max Id in
# Config to match.

$wb = $excel.ActiveWorkbook
$active = $wb.Worksheets.Item("Active")
$inactive = $wb.Worksheets.Item("Inactive")
$config = $wb.Worksheets.Item("Config")

# --- 1. Remove task #33 from Active (it was row 3) ---------------------
$active.Rows.Item(3).Delete()

# --- 2. Insert the two new Active tasks at their recorded positions ----
# #53 lands just above #51 (now row 10 after the delete above)
$active.Rows.Item(10).Insert()
$active.Cells.Item(10, 1).Value = 53
$active.Cells.Item(10, 2).Value = "prompt to save if image has changed since last save`n- on closing program`n- on opening new image"
$active.Cells.Item(10, 3).Value = "Todo"
$active.Cells.Item(10, 4).Value = "Task"
$active.Cells.Item(10, 5).Value = "'8/23/2018"

# #52 lands just above #27 (now row 13 after the insert above)
$active.Rows.Item(13).Insert()
$active.Cells.Item(13, 1).Value = 52
$active.Cells.Item(13, 2).Value = "in documentation`nrecommend user keeps an original b/w copy to go back to`nif conversion errors build up with lots of editing"
$active.Cells.Item(13, 3).Value = "Todo"
$active.Cells.Item(13, 4).Value = "Task"
$active.Cells.Item(13, 5).Value = "'8/22/2018"

# --- 3. Add task #33 to the top of Inactive, now Done ------------------
$inactive.Rows.Item(2).Insert()
$inactive.Cells.Item(2, 1).Value = 33
$inactive.Cells.Item(2, 2).Value = "undo, redo coloring a section on the image"
$inactive.Cells.Item(2, 3).Value = "Done"
$inactive.Cells.Item(2, 4).Value = "Task"
$inactive.Cells.Item(2, 5).Value = "'8/11/2018"
$inactive.Cells.Item(2, 6).Value = "'8/23/2018"

# --- 4. Bump the Max Id tracker in Config (51 -> 53) --------------------
$config.Cells.Item(2, 6).Value = 53
